# Regenerate save_data to use K instead of Strike#: update column G (K) values
# for rows 2-30 on Sheet1 with the newly-calculated strikeout (K) counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 6
    4  = 5
    5  = 5
    6  = 4
    7  = 3
    8  = 2
    9  = 1
    10 = 4
    11 = 1
    12 = 5
    13 = 2
    14 = 2
    15 = 3
    16 = 4
    17 = 4
    18 = 3
    19 = 2
    20 = 2
    21 = 4
    22 = 6
    23 = 4
    24 = 2
    25 = 2
    26 = 2
    27 = 7
    28 = 4
    29 = 4
    30 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
